$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Status and Date values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Re-apply the top/wrap-text cell alignment on both sheets so it is
#     explicitly flagged on the cell format (applyAlignment), matching the
#     published look of the implementation guide table. ---
$meta.Range("A1:B21").WrapText = $true
$meta.Range("A1:B21").VerticalAlignment = -4160

$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("A1:D6").WrapText = $true
$concepts.Range("A1:D6").VerticalAlignment = -4160
